$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '30.047.62'
$ws.Cells.Item(3, 4).Value = '1.872.44'
$ws.Cells.Item(5, 4).Value = '319.27'
$ws.Cells.Item(6, 4).Value = '1.001'
$ws.Cells.Item(7, 4).Value = '0.5049'
$ws.Cells.Item(8, 4).Value = '0.3964'
$ws.Cells.Item(9, 4).Value = '0.08213'
$ws.Cells.Item(10, 4).Value = '42.04'
$ws.Cells.Item(11, 4).Value = '1.094'
$ws.Cells.Item(12, 4).Value = '23.45'
$ws.Cells.Item(13, 4).Value = '1.879.73'
$ws.Cells.Item(14, 4).Value = '6.281'
$ws.Cells.Item(15, 4).Value = '7.197'
$ws.Cells.Item(17, 4).Value = '91.87'
$ws.Cells.Item(18, 4).Value = '0.00001086'
$ws.Cells.Item(19, 4).Value = '0.06451'
$ws.Cells.Item(20, 4).Value = '18.08'
$ws.Cells.Item(22, 4).Value = '30.034.84'
$ws.Cells.Item(23, 4).Value = '5.846'
$ws.Cells.Item(25, 4).Value = '2.170'
$ws.Cells.Item(26, 4).Value = '2.090.43'
$ws.Cells.Item(27, 4).Value = '21.32'
$ws.Cells.Item(28, 4).Value = '161.01'
$ws.Cells.Item(29, 4).Value = '2.226'
$ws.Cells.Item(30, 4).Value = '127.36'
$ws.Cells.Item(31, 4).Value = '1.075'
$ws.Cells.Item(33, 4).Value = '5.945'
$ws.Cells.Item(35, 4).Value = '0.02438'
$ws.Cells.Item(36, 4).Value = '5.237'
$ws.Cells.Item(37, 4).Value = '0.06369'
$ws.Cells.Item(39, 4).Value = '1.174'
$ws.Cells.Item(40, 4).Value = '8.483'
$ws.Cells.Item(41, 4).Value = '0.6306'
$ws.Cells.Item(42, 4).Value = '1.214'
$ws.Cells.Item(43, 4).Value = '11.31'
$ws.Cells.Item(44, 4).Value = '13.04'
$ws.Cells.Item(45, 4).Value = '0.5917'
$ws.Cells.Item(46, 4).Value = '2.111'
$ws.Cells.Item(48, 4).Value = '122.54'
$ws.Cells.Item(49, 4).Value = '1.207'
$ws.Cells.Item(50, 4).Value = '77.52'
$ws.Cells.Item(51, 4).Value = '1.116'

$ws.Cells.Item(2, 5).Value = '  -0.18%  '
$ws.Cells.Item(3, 5).Value = '  -2.39%  '
$ws.Cells.Item(5, 5).Value = '  -3.24%  '
$ws.Cells.Item(6, 5).Value = '  +0.05%  '
$ws.Cells.Item(7, 5).Value = '  -3.34%  '
$ws.Cells.Item(8, 5).Value = '  -3.06%  '
$ws.Cells.Item(9, 5).Value = '  -3.41%  '
$ws.Cells.Item(10, 5).Value = '  -2.22%  '
$ws.Cells.Item(11, 5).Value = '  -2.93%  '
$ws.Cells.Item(12, 5).Value = '  +4.18%  '
$ws.Cells.Item(13, 5).Value = '  -2.81%  '
$ws.Cells.Item(14, 5).Value = '  -2.16%  '
$ws.Cells.Item(15, 5).Value = '  -3.12%  '
$ws.Cells.Item(16, 5).Value = '  +0.08%  '
$ws.Cells.Item(17, 5).Value = '  -3.78%  '
$ws.Cells.Item(18, 5).Value = '  -2.36%  '
$ws.Cells.Item(19, 5).Value = '  -4.03%  '
$ws.Cells.Item(20, 5).Value = '  -1.33%  '
$ws.Cells.Item(21, 5).Value = '  +0.04%  '
$ws.Cells.Item(22, 5).Value = '  -0.22%  '
$ws.Cells.Item(23, 5).Value = '  -2.73%  '
$ws.Cells.Item(24, 5).Value = '  -1.69%  '
$ws.Cells.Item(25, 5).Value = '  -2.31%  '
$ws.Cells.Item(26, 5).Value = '  -2.94%  '
$ws.Cells.Item(27, 5).Value = '  +1.10%  '
$ws.Cells.Item(28, 5).Value = '  +0.36%  '
$ws.Cells.Item(29, 5).Value = '  -9.31%  '
$ws.Cells.Item(30, 5).Value = '  -1.53%  '
$ws.Cells.Item(31, 5).Value = '  -0.30%  '
$ws.Cells.Item(32, 5).Value = '  -2.02%  '
$ws.Cells.Item(33, 5).Value = '  -2.35%  '
$ws.Cells.Item(34, 5).Value = '  +1.47%  '
$ws.Cells.Item(35, 5).Value = '  -2.10%  '
$ws.Cells.Item(36, 5).Value = '  +0.87%  '
$ws.Cells.Item(37, 5).Value = '  -4.01%  '
$ws.Cells.Item(38, 5).Value = '  -3.16%  '
$ws.Cells.Item(39, 5).Value = '  -4.71%  '
$ws.Cells.Item(40, 5).Value = '  -4.80%  '
$ws.Cells.Item(41, 5).Value = '  -3.50%  '
$ws.Cells.Item(42, 5).Value = '  -2.64%  '
$ws.Cells.Item(43, 5).Value = '  -2.94%  '
$ws.Cells.Item(44, 5).Value = '  -2.03%  '
$ws.Cells.Item(45, 5).Value = '  -4.08%  '
$ws.Cells.Item(46, 5).Value = '  +1.25%  '
$ws.Cells.Item(47, 5).Value = '  -3.86%  '
$ws.Cells.Item(48, 5).Value = '  -1.67%  '
$ws.Cells.Item(49, 5).Value = '  -3.39%  '
$ws.Cells.Item(50, 5).Value = '  -2.95%  '
$ws.Cells.Item(51, 5).Value = '  -4.30%  '
